$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# BTec logo picture lives in the header that Word exposes as Headers.Item(2)
# (physically word/header1.xml); rename image1.jpg -> image2.jpg.
$hdr = $sec.Headers.Item(2)
if ($hdr.Range.InlineShapes.Count -gt 0) {
    $btecInline = $hdr.Range.InlineShapes.Item(1)
    $btecShape = $btecInline.ConvertToShape()
    $btecShape.Name = "image2.jpg"
    $btecShape.ConvertToInlineShape() | Out-Null
}

# Pearson logo picture #1 lives in the footer Word exposes as Footers.Item(1)
# (physically word/footer2.xml, docPr id="2"); rename image2.png -> image1.png.
$ftr1 = $sec.Footers.Item(1)
if ($ftr1.Range.InlineShapes.Count -gt 0) {
    $pearson1Inline = $ftr1.Range.InlineShapes.Item(1)
    $pearson1Shape = $pearson1Inline.ConvertToShape()
    $pearson1Shape.Name = "image1.png"
    $pearson1Shape.ConvertToInlineShape() | Out-Null
}

# Pearson logo picture #2 lives in the footer Word exposes as Footers.Item(2)
# (physically word/footer1.xml, docPr id="3"); rename image2.png -> image1.png.
$ftr2 = $sec.Footers.Item(2)
if ($ftr2.Range.InlineShapes.Count -gt 0) {
    $pearson2Inline = $ftr2.Range.InlineShapes.Item(1)
    $pearson2Shape = $pearson2Inline.ConvertToShape()
    $pearson2Shape.Name = "image1.png"
    $pearson2Shape.ConvertToInlineShape() | Out-Null
}

Write-Output "Renamed inline picture shapes."
